$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Collapse the three "CORE COMPETENCIES" bullet paragraphs into a
#    single summary paragraph.
# ---------------------------------------------------------------------
$coreHeading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("CORE COMPETENCIES")) {
        $coreHeading = $i
        break
    }
}

$p1 = $d.Paragraphs.Item($coreHeading + 1)
$p2 = $d.Paragraphs.Item($coreHeading + 2)
$p3 = $d.Paragraphs.Item($coreHeading + 3)

$bullet = [char]0x2022
$p1.Range.Text = "Product Marketing Core " + $bullet + " Research & Analytics " + $bullet + " Communication & Technology"

$deleteRange = $d.Range($p2.Range.Start, $p3.Range.End)
$deleteRange.Delete()

# ---------------------------------------------------------------------
# 2. Append a new "TECHNICAL SKILLS" section with the detailed bullet
#    lists that used to live under "CORE COMPETENCIES".
# ---------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastPara.Range.InsertParagraphAfter()

$lastIndex = $d.Paragraphs.Count
$headingPara = $d.Paragraphs.Item($lastIndex)
$headingPara.Range.Text = "TECHNICAL SKILLS"
$headingPara.Style = "Heading 2"

$headingPara.Range.InsertParagraphAfter()
$lastIndex = $d.Paragraphs.Count
$coreLine = $d.Paragraphs.Item($lastIndex)
$coreLine.Style = "Normal"
$coreLine.Range.Text = "PRODUCT MARKETING CORE Market Intelligence & Competitive Analysis; Product Positioning & Messaging Development; Go-to-Market Strategy & Product Launch Management; Customer Segmentation & Buyer Persona Development; Cross-functional Team Leadership & Collaboration; Sales Enablement & Training Material Development; Data-Driven Decision Making & Analytics Interpretation"

$coreLine.Range.InsertParagraphAfter()
$lastIndex = $d.Paragraphs.Count
$researchLine = $d.Paragraphs.Item($lastIndex)
$researchLine.Style = "Normal"
$researchLine.Range.Text = "RESEARCH & ANALYTICS Survey Methodology & Customer Insights; Market Research Design & Implementation; Competitive Intelligence & SWOT Analysis; Customer Journey Mapping & Behavioral Analysis; Statistical Modeling & Trend Analysis; Performance Metrics & Dashboard Development; A/B Testing & Conversion Optimization"

$researchLine.Range.InsertParagraphAfter()
$lastIndex = $d.Paragraphs.Count
$commLine = $d.Paragraphs.Item($lastIndex)
$commLine.Style = "Normal"
$commLine.Range.Text = "COMMUNICATION & TECHNOLOGY Strategic Messaging & Narrative Development; Technical Concept Translation for Business Audiences; Stakeholder Communication & Presentation Skills; Data Visualization & Reporting (Tableau, PowerBI, d3.js); Marketing Technology Stack Integration; Content Strategy & Thought Leadership; Client Relationship Management & Business Development"
